$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35 - this shifts existing rows 35..83 down to 36..84
$ws.Rows.Item(35).Insert()

# Populate the newly inserted row 35 with the new weekly price record
$ws.Cells.Item(35, 1).Value  = 4
$ws.Cells.Item(35, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(35, 3).Value  = "Los Lagos"
$ws.Cells.Item(35, 4).Value  = 44546
$ws.Cells.Item(35, 5).Value  = 10
$ws.Cells.Item(35, 6).Value  = 100112052
$ws.Cells.Item(35, 7).Value  = "Albahaca"
$ws.Cells.Item(35, 8).Value  = "Sin especificar"
$ws.Cells.Item(35, 9).Value  = "Primera"
$ws.Cells.Item(35, 10).Value = 60
$ws.Cells.Item(35, 11).Value = 8000
$ws.Cells.Item(35, 12).Value = 8000
$ws.Cells.Item(35, 13).Value = 8000
$ws.Cells.Item(35, 14).Value = "`$/docena de matas"
$ws.Cells.Item(35, 15).Value = "Región Metropolitana"
$ws.Cells.Item(35, 16).Value = 1333
$ws.Cells.Item(35, 17).Value = 6
$ws.Cells.Item(35, 18).Value = "Hortaliza"
